# Generate Report for Handback
#
# - Marks the "cc9c4d88-..." item as handed back (Status column) on both the
#   zh-cn and de-de report sheets.
# - Records the "Latest Target File" (F) / "Latest Handback File" (G)
#   hyperlinks for row 2 and row 3 on each sheet (new columns in the table).
# - Stamps the "Latest Handback DateTime" (H) for both rows on each sheet.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# Source (".md") file is the same for every row being reported on.
$mdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/bce50123aab168df4d4800b3a6c99b6650f81acf/e2e/cc9c4d88-a1e4-4904-acb4-021401f05c23.md"
$mdDisplay = "cc9c4d88-a1e4-4904-acb4-021401f05c23.md"

$sheetInfo = @(
    @{
        SheetName    = "zh-cn";
        XlfUrl       = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cbe94001b869f5cc9d6f2f6aab76a16704459c1d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/cc9c4d88-a1e4-4904-acb4-021401f05c23.6c0a929d12438973a48157a11a42e7268d01a887.zh-cn.xlf";
        XlfDisplay   = "cc9c4d88-a1e4-4904-acb4-021401f05c23.6c0a929d12438973a48157a11a42e7268d01a887.zh-cn.xlf";
        HandbackTime = "2016-03-17 14:51:10";
    },
    @{
        SheetName    = "de-de";
        XlfUrl       = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/601df4bf3caecfbd386fdd5f64d118538d3da328/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/cc9c4d88-a1e4-4904-acb4-021401f05c23.6c0a929d12438973a48157a11a42e7268d01a887.de-de.xlf";
        XlfDisplay   = "cc9c4d88-a1e4-4904-acb4-021401f05c23.6c0a929d12438973a48157a11a42e7268d01a887.de-de.xlf";
        HandbackTime = "2016-03-17 14:51:16";
    }
)

foreach ($info in $sheetInfo) {
    $ws = $wb.Worksheets.Item($info.SheetName)

    foreach ($rowNum in 2, 3) {
        # Status: the file is now handed back and in sync with en-US.
        $ws.Range("C" + $rowNum).Value = $statusHandedBack

        # Latest Target File (F) / Latest Handback File (G): add the
        # handback report columns with hyperlinks to the source .md and the
        # generated .xlf handback file.
        $ws.Hyperlinks.Add($ws.Range("F" + $rowNum), $mdUrl, "", "", $mdDisplay)
        $ws.Hyperlinks.Add($ws.Range("G" + $rowNum), $info.XlfUrl, "", "", $info.XlfDisplay)

        # Latest Handback DateTime (H): stamp with the handback time.
        $ws.Range("H" + $rowNum).Value = $info.HandbackTime
    }
}
